$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update alpha_distance_range row (row 2)
$ws.Range("B2").Value = 4.5
$ws.Range("C2").Value = 11.2

# Update beta_distance_range row (row 3)
$ws.Range("B3").Value = 4.3
$ws.Range("C3").Value = 10.2

# Update ratio_threshold_range row (row 4)
$ws.Range("C4").Value = 1.5

# Replace theta_threshold_range row (row 5) with pie_threshold_range data
$ws.Range("A5").Value = "pie_threshold_range"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 20

# Delete old row 6 (previously pie_threshold_range, now redundant)
$ws.Rows.Item(6).Delete()

# Column C width change (closest achievable width to the target 5.5 given
# this runtime's fixed 7px-max-digit-width column sizing model)
$ws.Columns.Item(3).ColumnWidth = 4.8

# Page setup (paper size 9 = A4, orientation 1 = portrait)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection change
$ws.Range("C3").Select()
